# Rename the "_old"/"_new" suffixed header labels to the format-version
# specific suffixes "_FV2404"/"_FV2410", then turn the data range into a
# proper Excel Table and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A-J hold the "_old" (-> "_FV2404") headers, columns L-U hold the
# "_new" (-> "_FV2410") headers; column K ("diff") is left untouched.
for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $oldCol = $i + 1
    $newCol = $i + 12
    $ws.Cells.Item(1, $oldCol).Value = $baseNames[$i] + "_FV2404"
    $ws.Cells.Item(1, $newCol).Value = $baseNames[$i] + "_FV2410"
}

# Turn the used range into a proper table.
$rng = $ws.Range("A1:U59")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"

# Freeze the header row.
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
